$d = $word.ActiveDocument

# The Pearson/BTec logos living in this document's headers and footers were
# exported with their non-visual drawing names swapped relative to what they
# should be (Pearson logo tagged "image1.png" instead of "image2.png", BTec
# logo tagged "image2.jpg" instead of "image1.jpg"). Walk every header and
# footer in every section, find the inline picture(s), and rename each one
# according to that rule.

function Rename-InlineLogo($inlineShape) {
    if ($inlineShape -eq $null) { return }

    $shape = $inlineShape.ConvertToShape()
    $oldName = $shape.Name

    $newName = $null
    if ($oldName -eq "image1.png") {
        $newName = "image2.png"
    } elseif ($oldName -eq "image2.jpg") {
        $newName = "image1.jpg"
    }

    if ($newName -ne $null) {
        $shape.Name = $newName
    }

    [void]$shape.ConvertToInlineShape()
}

foreach ($sec in $d.Sections) {
    for ($i = 1; $i -le 3; $i++) {
        $hdr = $sec.Headers.Item($i)
        if ($hdr.Exists) {
            $shapes = $hdr.Range.InlineShapes
            for ($j = 1; $j -le $shapes.Count; $j++) {
                Rename-InlineLogo $shapes.Item($j)
            }
        }

        $ftr = $sec.Footers.Item($i)
        if ($ftr.Exists) {
            $shapes = $ftr.Range.InlineShapes
            for ($j = 1; $j -le $shapes.Count; $j++) {
                Rename-InlineLogo $shapes.Item($j)
            }
        }
    }
}

Write-Host "Renamed logo inline shapes."
